$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.338.19'
$ws.Range("E2").Value = '  +2.20%  '

$ws.Range("D3").Value = '1.661.80'
$ws.Range("E3").Value = '  +1.24%  '

$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").Value = '220.26'
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("E8").Value = '  +1.46%  '

$ws.Range("D10").Value = '19.96'
$ws.Range("E10").Value = '  +4.18%  '

$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").Value = '1.894.11'
$ws.Range("E12").Value = '  +1.27%  '

$ws.Range("D13").Value = '1.662.27'
$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("E15").Value = '  +1.57%  '

$ws.Range("D16").Value = '67.47'
$ws.Range("E16").Value = '  +4.28%  '

$ws.Range("D17").Value = '27.322.90'
$ws.Range("E17").Value = '  +2.09%  '

$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").Value = '222.96'
$ws.Range("E19").Value = '  +3.57%  '

$ws.Range("E20").Value = '  -0.27%  '

$ws.Range("D21").Value = '6.78'
$ws.Range("E21").Value = '  +8.74%  '

$ws.Range("E22").Value = '  +1.80%  '

$ws.Range("D23").Value = '2.50'
$ws.Range("E23").Value = '  +4.68%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").Value = '147.43'
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("E27").Value = '  +3.70%  '

$ws.Range("D28").Value = '0.119'
$ws.Range("E28").Value = '  +0.78%  '

$ws.Range("D29").Value = '16.05'
$ws.Range("E29").Value = '  +2.60%  '

$ws.Range("D30").Value = '0.0516'
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("E31").Value = '  +0.93%  '

$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("E33").Value = '  +0.44%  '

$ws.Range("E34").Value = '  +2.02%  '

$ws.Range("D35").Value = '1.260.06'
$ws.Range("E35").Value = '  -1.96%  '

$ws.Range("E36").Value = '  +0.27%  '

$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("E39").Value = '  +1.88%  '

$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("E42").Value = '  +2.07%  '

$ws.Range("D43").Value = '1.806.07'
$ws.Range("E43").Value = '  +1.52%  '

$ws.Range("E44").Value = '  -4.15%  '

$ws.Range("D45").Value = '61.94'
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").Value = '92.59'
$ws.Range("E46").Value = '  +0.65%  '

$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.0986'
$ws.Range("E49").Value = '  +1.96%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.69'
$ws.Range("E50").Value = '  +0.45%  '

$ws.Range("E51").Value = '  +0.29%  '
